$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 231.36363
$ws.Cells.Item(2, 9).Value = 231.36363
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 231.36363
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -118.36363
$ws.Cells.Item(2, 14).ClearContents()
$ws.Cells.Item(12, 8).Value = 500
$ws.Cells.Item(12, 10).Value = 700
$ws.Cells.Item(12, 12).Value = 700
$ws.Cells.Item(12, 14).Value = -1040
$ws.Cells.Item(38, 8).Value = 2732.2
$ws.Cells.Item(38, 9).Value = 553.6667
$ws.Cells.Item(38, 10).Value = 6000
$ws.Cells.Item(38, 11).Value = 1661.0001
$ws.Cells.Item(38, 12).Value = 18000
$ws.Cells.Item(38, 13).Value = -1289.0001
$ws.Cells.Item(38, 14).Value = -18744
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 13).ClearContents()
$ws.Cells.Item(60, 8).Value = 0
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 13).ClearContents()
$ws.Cells.Item(64, 8).Value = 3499.3333
$ws.Cells.Item(64, 9).Value = 3499
$ws.Cells.Item(64, 11).Value = 3499
$ws.Cells.Item(64, 13).Value = -3251
$ws.Cells.Item(67, 8).Value = 3499.3333
$ws.Cells.Item(67, 9).Value = 3499
$ws.Cells.Item(67, 11).Value = 3499
$ws.Cells.Item(67, 13).Value = -2641
$ws.Cells.Item(82, 8).Value = 6588
$ws.Cells.Item(82, 9).Value = 4029.75
$ws.Cells.Item(82, 11).Value = 12089.25
$ws.Cells.Item(82, 13).Value = -11683.25
$ws.Cells.Item(85, 8).Value = 6588
$ws.Cells.Item(85, 9).Value = 4029.75
$ws.Cells.Item(85, 11).Value = 12089.25
$ws.Cells.Item(85, 13).Value = -10685.25
$ws.Cells.Item(92, 8).Value = 744.2143
$ws.Cells.Item(92, 9).Value = 744.2143
$ws.Cells.Item(92, 11).Value = 744.2143
$ws.Cells.Item(92, 13).Value = 503.7857
$ws.Cells.Item(100, 8).Value = 5127
$ws.Cells.Item(100, 9).Value = 5710.778
$ws.Cells.Item(100, 11).Value = 5710.778
$ws.Cells.Item(100, 13).Value = -5169.778
$ws.Cells.Item(103, 8).Value = 4187.875
$ws.Cells.Item(103, 9).Value = 4571.857
$ws.Cells.Item(103, 11).Value = 13715.571
$ws.Cells.Item(103, 13).Value = -13129.571
$ws.Cells.Item(112, 8).Value = 2092.8572
$ws.Cells.Item(112, 10).Value = 2185.4614
$ws.Cells.Item(112, 12).Value = 6556.3842
$ws.Cells.Item(112, 14).Value = -8772.3842
$ws.Cells.Item(121, 8).Value = 1838.3846
$ws.Cells.Item(121, 10).Value = 1873.091
$ws.Cells.Item(121, 12).Value = 5619.272999999999
$ws.Cells.Item(121, 14).Value = -9113.272999999999
$ws.Cells.Item(123, 8).Value = 109985.4
$ws.Cells.Item(123, 10).Value = 109985.4
$ws.Cells.Item(123, 12).Value = 109985.4
$ws.Cells.Item(123, 14).Value = -119785.4
$ws.Cells.Item(124, 8).Value = 193383.33
$ws.Cells.Item(124, 10).Value = 193383.33
$ws.Cells.Item(124, 12).Value = 193383.33
$ws.Cells.Item(124, 14).Value = -203203.33
$ws.Cells.Item(129, 8).Value = 2088.3333
$ws.Cells.Item(129, 9).Value = 1382.75
$ws.Cells.Item(129, 11).Value = 4148.25
$ws.Cells.Item(129, 13).Value = 851.75
$ws.Cells.Item(132, 8).Value = 5109.8125
$ws.Cells.Item(132, 9).Value = 5204.2964
$ws.Cells.Item(132, 11).Value = 15612.8892
$ws.Cells.Item(132, 13).Value = -13082.8892
$ws.Cells.Item(134, 8).Value = 65226.773
$ws.Cells.Item(134, 9).Value = 20000
$ws.Cells.Item(134, 11).Value = 20000
$ws.Cells.Item(134, 13).Value = -14930
$ws.Cells.Item(135, 8).Value = 1475.2778
$ws.Cells.Item(135, 9).Value = 1337.0667
$ws.Cells.Item(135, 11).Value = 12033.6003
$ws.Cells.Item(135, 13).Value = -9498.6003
$ws.Cells.Item(137, 8).Value = 1858335.5
$ws.Cells.Item(137, 9).Value = 2273843.2
$ws.Cells.Item(137, 10).Value = 30101.6
$ws.Cells.Item(137, 11).Value = 6821529.600000001
$ws.Cells.Item(137, 12).Value = 90304.79999999999
$ws.Cells.Item(137, 13).Value = -6818979.600000001
$ws.Cells.Item(137, 14).Value = -95404.79999999999
$ws.Cells.Item(138, 8).Value = 6110.087
$ws.Cells.Item(138, 9).Value = 3741.625
$ws.Cells.Item(138, 10).Value = 7373.2666
$ws.Cells.Item(138, 11).Value = 11224.875
$ws.Cells.Item(138, 12).Value = 22119.7998
$ws.Cells.Item(138, 13).Value = -6084.875
$ws.Cells.Item(138, 14).Value = -32399.7998
$ws.Cells.Item(139, 8).Value = 69970.71000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1863.5405
$ws.Cells.Item(2, 9).Value = 805.0968
$ws.Cells.Item(2, 11).Value = 805.0968
$ws.Cells.Item(2, 13).Value = -692.0968
$ws.Cells.Item(5, 8).Value = 66.375
$ws.Cells.Item(32, 8).Value = 2568554.2
$ws.Cells.Item(32, 9).Value = 1326340.4
$ws.Cells.Item(32, 10).Value = 9524952
$ws.Cells.Item(32, 11).Value = 1326340.4
$ws.Cells.Item(32, 12).Value = 9524952
$ws.Cells.Item(32, 13).Value = -1326053.4
$ws.Cells.Item(32, 14).Value = -9525526
$ws.Cells.Item(45, 8).Value = 71537350
$ws.Cells.Item(45, 9).Value = 149288.4
$ws.Cells.Item(45, 11).Value = 149288.4
$ws.Cells.Item(45, 13).Value = -148911.4
$ws.Cells.Item(61, 8).Value = 3665.2
$ws.Cells.Item(61, 9).Value = 4163.1665
$ws.Cells.Item(61, 10).Value = 3333.2222
$ws.Cells.Item(61, 11).Value = 4163.1665
$ws.Cells.Item(61, 12).Value = 3333.2222
$ws.Cells.Item(61, 13).Value = -3951.1665
$ws.Cells.Item(61, 14).Value = -3757.2222
$ws.Cells.Item(74, 8).Value = 15237728
$ws.Cells.Item(74, 9).Value = 130703.35
$ws.Cells.Item(74, 10).Value = 43481296
$ws.Cells.Item(74, 11).Value = 130703.35
$ws.Cells.Item(74, 12).Value = 43481296
$ws.Cells.Item(74, 13).Value = -129829.35
$ws.Cells.Item(74, 14).Value = -43483044
$ws.Cells.Item(77, 8).Value = 15237728
$ws.Cells.Item(77, 9).Value = 130703.35
$ws.Cells.Item(77, 10).Value = 43481296
$ws.Cells.Item(77, 11).Value = 653516.75
$ws.Cells.Item(77, 12).Value = 217406480
$ws.Cells.Item(77, 13).Value = -649148.75
$ws.Cells.Item(77, 14).Value = -217415216
$ws.Cells.Item(97, 8).Value = 1760.8823
$ws.Cells.Item(97, 9).Value = 1760.8823
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 1760.8823
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = -1264.8823
$ws.Cells.Item(97, 14).ClearContents()
$ws.Cells.Item(102, 8).Value = 2169.2
$ws.Cells.Item(102, 9).Value = 1650
$ws.Cells.Item(102, 11).Value = 1650
$ws.Cells.Item(102, 13).Value = -28
$ws.Cells.Item(110, 8).Value = 1803.2609
$ws.Cells.Item(110, 9).Value = 1784.5714
$ws.Cells.Item(110, 10).Value = 1999.5
$ws.Cells.Item(110, 11).Value = 1784.5714
$ws.Cells.Item(110, 12).Value = 1999.5
$ws.Cells.Item(110, 13).Value = 260.4286
$ws.Cells.Item(110, 14).Value = -6089.5
$ws.Cells.Item(116, 8).Value = 1863.5405
$ws.Cells.Item(116, 9).Value = 805.0968
$ws.Cells.Item(116, 11).Value = 805.0968
$ws.Cells.Item(116, 13).Value = 1488.9032
$ws.Cells.Item(122, 8).Value = 3065.4092
$ws.Cells.Item(122, 9).Value = 2445.7856
$ws.Cells.Item(122, 10).Value = 4149.75
$ws.Cells.Item(122, 11).Value = 7337.3568
$ws.Cells.Item(122, 12).Value = 12449.25
$ws.Cells.Item(122, 13).Value = -4887.3568
$ws.Cells.Item(122, 14).Value = -17349.25
$ws.Cells.Item(132, 8).Value = 1510158.2
$ws.Cells.Item(132, 9).Value = 2025866.4
$ws.Cells.Item(132, 10).Value = 2704.077
$ws.Cells.Item(132, 11).Value = 6077599.199999999
$ws.Cells.Item(132, 12).Value = 8112.231000000001
$ws.Cells.Item(132, 13).Value = -6075069.199999999
$ws.Cells.Item(132, 14).Value = -13172.231
$ws.Cells.Item(136, 8).Value = 3665.2
$ws.Cells.Item(136, 9).Value = 4163.1665
$ws.Cells.Item(136, 10).Value = 3333.2222
$ws.Cells.Item(136, 11).Value = 12489.4995
$ws.Cells.Item(136, 12).Value = 9999.6666
$ws.Cells.Item(136, 13).Value = -9939.499500000002
$ws.Cells.Item(136, 14).Value = -15099.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1863.5405
$ws.Cells.Item(3, 9).Value = 805.0968
$ws.Cells.Item(3, 11).Value = 805.0968
$ws.Cells.Item(3, 13).Value = -691.0968
$ws.Cells.Item(4, 8).Value = 66.375
$ws.Cells.Item(20, 8).Value = 24514180
$ws.Cells.Item(20, 9).Value = 29766892
$ws.Cells.Item(20, 11).Value = 29766892
$ws.Cells.Item(20, 13).Value = -29766645
$ws.Cells.Item(22, 8).Value = 292.08334
$ws.Cells.Item(22, 9).Value = 320.8889
$ws.Cells.Item(22, 10).Value = 205.66667
$ws.Cells.Item(22, 11).Value = 320.8889
$ws.Cells.Item(22, 12).Value = 205.66667
$ws.Cells.Item(22, 13).Value = -147.8889
$ws.Cells.Item(22, 14).Value = -551.6666700000001
$ws.Cells.Item(80, 8).Value = 939.6
$ws.Cells.Item(80, 10).Value = 948.75
$ws.Cells.Item(80, 12).Value = 948.75
$ws.Cells.Item(80, 14).Value = -2944.75
$ws.Cells.Item(81, 8).Value = 16355.625
$ws.Cells.Item(81, 10).Value = 16355.625
$ws.Cells.Item(81, 12).Value = 16355.625
$ws.Cells.Item(81, 14).Value = -18477.625
$ws.Cells.Item(83, 8).Value = 939.6
$ws.Cells.Item(83, 10).Value = 948.75
$ws.Cells.Item(83, 12).Value = 4743.75
$ws.Cells.Item(83, 14).Value = -14727.75
$ws.Cells.Item(84, 8).Value = 16355.625
$ws.Cells.Item(84, 10).Value = 16355.625
$ws.Cells.Item(84, 12).Value = 49066.875
$ws.Cells.Item(84, 14).Value = -59674.875
$ws.Cells.Item(94, 8).Value = 44445670
$ws.Cells.Item(94, 9).Value = 47620332
$ws.Cells.Item(94, 10).Value = 472.5
$ws.Cells.Item(94, 11).Value = 47620332
$ws.Cells.Item(94, 12).Value = 472.5
$ws.Cells.Item(94, 13).Value = -47619881
$ws.Cells.Item(94, 14).Value = -1374.5
$ws.Cells.Item(99, 8).Value = 5765
$ws.Cells.Item(99, 9).Value = 5418
$ws.Cells.Item(99, 11).Value = 5418
$ws.Cells.Item(99, 13).Value = -3920
$ws.Cells.Item(105, 8).Value = 15296803
$ws.Cells.Item(105, 10).Value = 27781130
$ws.Cells.Item(105, 12).Value = 27781130
$ws.Cells.Item(105, 14).Value = -27784624
$ws.Cells.Item(107, 8).Value = 1973865
$ws.Cells.Item(107, 9).Value = 2332356.2
$ws.Cells.Item(107, 11).Value = 2332356.2
$ws.Cells.Item(107, 13).Value = -2330436.2
$ws.Cells.Item(132, 8).Value = 89799.60000000001
$ws.Cells.Item(132, 10).Value = 89799.60000000001
$ws.Cells.Item(132, 12).Value = 89799.60000000001
$ws.Cells.Item(132, 14).Value = -99919.60000000001
$ws.Cells.Item(134, 8).Value = 2541.1667
$ws.Cells.Item(134, 9).Value = 1780.875
$ws.Cells.Item(134, 10).Value = 3149.4
$ws.Cells.Item(134, 11).Value = 5342.625
$ws.Cells.Item(134, 12).Value = 9448.200000000001
$ws.Cells.Item(134, 13).Value = -2807.625
$ws.Cells.Item(134, 14).Value = -14518.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 831.875
$ws.Cells.Item(16, 9).Value = 848.4
$ws.Cells.Item(16, 11).Value = 848.4
$ws.Cells.Item(16, 13).Value = -561.4
$ws.Cells.Item(31, 8).Value = 1669648.1
$ws.Cells.Item(31, 9).Value = 1666.3636
$ws.Cells.Item(31, 10).Value = 1956332.5
$ws.Cells.Item(31, 11).Value = 1666.3636
$ws.Cells.Item(31, 12).Value = 1956332.5
$ws.Cells.Item(31, 13).Value = -1371.3636
$ws.Cells.Item(31, 14).Value = -1956922.5
$ws.Cells.Item(34, 8).Value = 1669648.1
$ws.Cells.Item(34, 9).Value = 1666.3636
$ws.Cells.Item(34, 10).Value = 1956332.5
$ws.Cells.Item(34, 11).Value = 1666.3636
$ws.Cells.Item(34, 12).Value = 1956332.5
$ws.Cells.Item(34, 13).Value = -1464.3636
$ws.Cells.Item(34, 14).Value = -1956736.5
$ws.Cells.Item(53, 8).Value = 78594.75
$ws.Cells.Item(53, 10).Value = 78594.75
$ws.Cells.Item(53, 12).Value = 78594.75
$ws.Cells.Item(53, 14).Value = -79808.75
$ws.Cells.Item(58, 8).Value = 5620.913
$ws.Cells.Item(58, 9).Value = 3708
$ws.Cells.Item(58, 10).Value = 7092.385
$ws.Cells.Item(58, 11).Value = 3708
$ws.Cells.Item(58, 12).Value = 7092.385
$ws.Cells.Item(58, 13).Value = -3505
$ws.Cells.Item(58, 14).Value = -7498.385
$ws.Cells.Item(62, 8).Value = 3541.0557
$ws.Cells.Item(62, 9).Value = 3615.75
$ws.Cells.Item(62, 11).Value = 3615.75
$ws.Cells.Item(62, 13).Value = -2991.75
$ws.Cells.Item(65, 8).Value = 3541.0557
$ws.Cells.Item(65, 9).Value = 3615.75
$ws.Cells.Item(65, 11).Value = 18078.75
$ws.Cells.Item(65, 13).Value = -14958.75
$ws.Cells.Item(94, 8).Value = 720.5
$ws.Cells.Item(94, 9).Value = 414.5
$ws.Cells.Item(94, 10).Value = 873.5
$ws.Cells.Item(94, 11).Value = 414.5
$ws.Cells.Item(94, 12).Value = 873.5
$ws.Cells.Item(94, 13).Value = 36.5
$ws.Cells.Item(94, 14).Value = -1775.5
$ws.Cells.Item(113, 8).Value = 831.875
$ws.Cells.Item(113, 9).Value = 848.4
$ws.Cells.Item(113, 11).Value = 848.4
$ws.Cells.Item(113, 13).Value = 1321.6
$ws.Cells.Item(132, 8).Value = 2887.024
$ws.Cells.Item(132, 9).Value = 2627.4375
$ws.Cells.Item(132, 11).Value = 7882.3125
$ws.Cells.Item(132, 13).Value = -5352.3125
$ws.Cells.Item(134, 8).Value = 3610.0344
$ws.Cells.Item(134, 9).Value = 3462.1428
$ws.Cells.Item(134, 10).Value = 3998.25
$ws.Cells.Item(134, 11).Value = 10386.4284
$ws.Cells.Item(134, 12).Value = 11994.75
$ws.Cells.Item(134, 13).Value = -7851.428400000001
$ws.Cells.Item(134, 14).Value = -17064.75
$ws.Cells.Item(136, 8).Value = 5620.913
$ws.Cells.Item(136, 9).Value = 3708
$ws.Cells.Item(136, 10).Value = 7092.385
$ws.Cells.Item(136, 11).Value = 11124
$ws.Cells.Item(136, 12).Value = 21277.155
$ws.Cells.Item(136, 13).Value = -8574
$ws.Cells.Item(136, 14).Value = -26377.155
$ws.Cells.Item(138, 8).Value = 69998.75
$ws.Cells.Item(138, 10).Value = 69998.75
$ws.Cells.Item(138, 12).Value = 69998.75
$ws.Cells.Item(138, 14).Value = -80278.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1908749.6
$ws.Cells.Item(4, 10).Value = 8242204
$ws.Cells.Item(4, 12).Value = 24726612
$ws.Cells.Item(4, 14).Value = -24726836
$ws.Cells.Item(5, 8).Value = 419.77777
$ws.Cells.Item(5, 9).Value = 391.16666
$ws.Cells.Item(5, 10).Value = 477
$ws.Cells.Item(5, 11).Value = 1173.49998
$ws.Cells.Item(5, 12).Value = 1431
$ws.Cells.Item(5, 13).Value = -1061.49998
$ws.Cells.Item(5, 14).Value = -1655
$ws.Cells.Item(11, 8).Value = 111849.336
$ws.Cells.Item(11, 9).Value = 111849.336
$ws.Cells.Item(11, 11).Value = 335548.008
$ws.Cells.Item(11, 13).Value = -335408.008
$ws.Cells.Item(26, 8).Value = 1449.4117
$ws.Cells.Item(26, 9).Value = 1449.4117
$ws.Cells.Item(26, 11).Value = 4348.2351
$ws.Cells.Item(26, 13).Value = -4060.2351
$ws.Cells.Item(37, 8).Value = 98181.82000000001
$ws.Cells.Item(37, 10).Value = 98181.82000000001
$ws.Cells.Item(37, 12).Value = 294545.46
$ws.Cells.Item(37, 14).Value = -294769.46
$ws.Cells.Item(38, 8).Value = 397.9091
$ws.Cells.Item(38, 10).Value = 452.77777
$ws.Cells.Item(38, 12).Value = 1358.33331
$ws.Cells.Item(38, 14).Value = -2052.33331
$ws.Cells.Item(50, 8).Value = 222
$ws.Cells.Item(50, 10).Value = 222
$ws.Cells.Item(50, 12).Value = 666
$ws.Cells.Item(50, 14).Value = -1628
$ws.Cells.Item(53, 8).Value = 222
$ws.Cells.Item(53, 10).Value = 222
$ws.Cells.Item(53, 12).Value = 666
$ws.Cells.Item(53, 14).Value = -1628
$ws.Cells.Item(56, 8).Value = 7719.5884
$ws.Cells.Item(56, 9).Value = 7719.5884
$ws.Cells.Item(56, 11).Value = 7719.5884
$ws.Cells.Item(56, 13).Value = -7189.5884
$ws.Cells.Item(68, 8).Value = 6064849.5
$ws.Cells.Item(68, 10).Value = 6255022.5
$ws.Cells.Item(68, 12).Value = 18765067.5
$ws.Cells.Item(68, 14).Value = -18766689.5
$ws.Cells.Item(71, 8).Value = 6064849.5
$ws.Cells.Item(71, 10).Value = 6255022.5
$ws.Cells.Item(71, 12).Value = 56295202.5
$ws.Cells.Item(71, 14).Value = -56303314.5
$ws.Cells.Item(107, 8).Value = 2881.5
$ws.Cells.Item(107, 10).Value = 6862.6
$ws.Cells.Item(107, 12).Value = 20587.8
$ws.Cells.Item(107, 14).Value = -24427.8
$ws.Cells.Item(113, 8).Value = 805.8125
$ws.Cells.Item(113, 10).Value = 902.6667
$ws.Cells.Item(113, 12).Value = 2708.0001
$ws.Cells.Item(113, 14).Value = -7048.0001
$ws.Cells.Item(122, 8).Value = 1763.24
$ws.Cells.Item(122, 9).Value = 1257.5555
$ws.Cells.Item(122, 10).Value = 2047.6875
$ws.Cells.Item(122, 11).Value = 11317.9995
$ws.Cells.Item(122, 12).Value = 18429.1875
$ws.Cells.Item(122, 13).Value = -8867.9995
$ws.Cells.Item(122, 14).Value = -23329.1875
$ws.Cells.Item(128, 8).Value = 111992.25
$ws.Cells.Item(128, 9).Value = 111992.25
$ws.Cells.Item(128, 11).Value = 335976.75
$ws.Cells.Item(128, 13).Value = -330996.75
$ws.Cells.Item(131, 8).Value = 3814.2258
$ws.Cells.Item(131, 9).Value = 12565.6
$ws.Cells.Item(131, 11).Value = 37696.8
$ws.Cells.Item(131, 13).Value = -32656.8
$ws.Cells.Item(135, 8).Value = 419.77777
$ws.Cells.Item(135, 9).Value = 391.16666
$ws.Cells.Item(135, 10).Value = 477
$ws.Cells.Item(135, 11).Value = 3520.49994
$ws.Cells.Item(135, 12).Value = 4293
$ws.Cells.Item(135, 13).Value = -985.4999399999997
$ws.Cells.Item(135, 14).Value = -9363
$ws.Cells.Item(137, 8).Value = 6872.6
$ws.Cells.Item(137, 10).Value = 2666.3333
$ws.Cells.Item(137, 12).Value = 7998.999899999999
$ws.Cells.Item(137, 14).Value = -18198.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 66670750
$ws.Cells.Item(80, 10).Value = 5749.75
$ws.Cells.Item(80, 12).Value = 5749.75
$ws.Cells.Item(80, 14).Value = -7745.75
$ws.Cells.Item(83, 8).Value = 66670750
$ws.Cells.Item(83, 10).Value = 5749.75
$ws.Cells.Item(83, 12).Value = 28748.75
$ws.Cells.Item(83, 14).Value = -38732.75
$ws.Cells.Item(107, 8).Value = 20463.334
$ws.Cells.Item(107, 9).Value = 445.5
$ws.Cells.Item(107, 10).Value = 60499
$ws.Cells.Item(107, 11).Value = 445.5
$ws.Cells.Item(107, 12).Value = 60499
$ws.Cells.Item(107, 13).Value = 1474.5
$ws.Cells.Item(107, 14).Value = -64339
$ws.Cells.Item(113, 8).Value = 200004300
$ws.Cells.Item(113, 9).Value = 250004130
$ws.Cells.Item(113, 10).Value = 5000
$ws.Cells.Item(113, 11).Value = 250004130
$ws.Cells.Item(113, 12).Value = 5000
$ws.Cells.Item(113, 13).Value = -250001960
$ws.Cells.Item(113, 14).Value = -9340
$ws.Cells.Item(132, 8).Value = 2319.342
$ws.Cells.Item(132, 9).Value = 2349.5833
$ws.Cells.Item(132, 10).Value = 2267.5
$ws.Cells.Item(132, 11).Value = 7048.749899999999
$ws.Cells.Item(132, 12).Value = 6802.5
$ws.Cells.Item(132, 13).Value = -4518.749899999999
$ws.Cells.Item(132, 14).Value = -11862.5
$ws.Cells.Item(140, 8).Value = 69768

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 14).ClearContents()
$ws.Cells.Item(16, 8).Value = 3484.75
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 13).ClearContents()
$ws.Cells.Item(22, 8).Value = 333333730
$ws.Cells.Item(22, 9).Value = 333333730
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 333333730
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -333333435
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(27, 8).Value = 333333730
$ws.Cells.Item(27, 9).Value = 333333730
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 333333730
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = -333333623
$ws.Cells.Item(27, 14).ClearContents()
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 14).ClearContents()
$ws.Cells.Item(46, 8).Value = 956.25
$ws.Cells.Item(46, 9).Value = 878.5714
$ws.Cells.Item(46, 11).Value = 878.5714
$ws.Cells.Item(46, 13).Value = -690.5714
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 14).ClearContents()
$ws.Cells.Item(51, 8).Value = 30076
$ws.Cells.Item(51, 9).Value = 30076
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 30076
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = -29598
$ws.Cells.Item(51, 14).ClearContents()
$ws.Cells.Item(56, 8).Value = 18762.25
$ws.Cells.Item(56, 9).Value = 51
$ws.Cells.Item(56, 11).Value = 51
$ws.Cells.Item(56, 13).Value = 640
$ws.Cells.Item(61, 8).Value = 6075.923
$ws.Cells.Item(61, 9).Value = 6078.3
$ws.Cells.Item(61, 10).Value = 6068
$ws.Cells.Item(61, 11).Value = 6078.3
$ws.Cells.Item(61, 12).Value = 6068
$ws.Cells.Item(61, 13).Value = -5876.3
$ws.Cells.Item(61, 14).Value = -6472
$ws.Cells.Item(68, 8).Value = 33336534
$ws.Cells.Item(68, 9).Value = 47622336
$ws.Cells.Item(68, 11).Value = 47622336
$ws.Cells.Item(68, 13).Value = -47621587
$ws.Cells.Item(71, 8).Value = 33336534
$ws.Cells.Item(71, 9).Value = 47622336
$ws.Cells.Item(71, 11).Value = 238111680
$ws.Cells.Item(71, 13).Value = -238107936
$ws.Cells.Item(100, 8).Value = 3187.5
$ws.Cells.Item(100, 9).Value = 3187.5
$ws.Cells.Item(100, 11).Value = 3187.5
$ws.Cells.Item(100, 13).Value = -2646.5
$ws.Cells.Item(113, 8).Value = 6075.923
$ws.Cells.Item(113, 9).Value = 6078.3
$ws.Cells.Item(113, 10).Value = 6068
$ws.Cells.Item(113, 11).Value = 6078.3
$ws.Cells.Item(113, 12).Value = 6068
$ws.Cells.Item(113, 13).Value = -3908.3
$ws.Cells.Item(113, 14).Value = -10408
$ws.Cells.Item(132, 8).Value = 5933.3335
$ws.Cells.Item(132, 9).Value = 6939.9
$ws.Cells.Item(132, 10).Value = 5018.273
$ws.Cells.Item(132, 11).Value = 20819.7
$ws.Cells.Item(132, 12).Value = 15054.819
$ws.Cells.Item(132, 13).Value = -18289.7
$ws.Cells.Item(132, 14).Value = -20114.819
$ws.Cells.Item(133, 8).Value = 107775
$ws.Cells.Item(133, 10).Value = 107775
$ws.Cells.Item(133, 12).Value = 107775
$ws.Cells.Item(133, 14).Value = -112835
$ws.Cells.Item(134, 8).Value = 109496.5
$ws.Cells.Item(134, 10).Value = 109496.5
$ws.Cells.Item(134, 12).Value = 109496.5
$ws.Cells.Item(134, 14).Value = -119636.5
$ws.Cells.Item(135, 8).Value = 100000
$ws.Cells.Item(135, 9).Value = 100000
$ws.Cells.Item(135, 11).Value = 100000
$ws.Cells.Item(135, 13).Value = -94930
$ws.Cells.Item(136, 8).Value = 5797.6553
$ws.Cells.Item(136, 9).Value = 4087.5293
$ws.Cells.Item(136, 10).Value = 8220.333000000001
$ws.Cells.Item(136, 11).Value = 12262.5879
$ws.Cells.Item(136, 12).Value = 24660.999
$ws.Cells.Item(136, 13).Value = -9712.5879
$ws.Cells.Item(136, 14).Value = -29760.999
$ws.Cells.Item(138, 8).Value = 92900
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()
$ws.Cells.Item(139, 8).Value = 70040.45
$ws.Cells.Item(139, 10).Value = 70040.45
$ws.Cells.Item(139, 12).Value = 70040.45
$ws.Cells.Item(139, 14).Value = -80320.45
$ws.Cells.Item(140, 8).Value = 96790
$ws.Cells.Item(140, 10).Value = 96790
$ws.Cells.Item(140, 12).Value = 96790
$ws.Cells.Item(140, 14).Value = -107150

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 13).ClearContents()
$ws.Cells.Item(46, 8).Value = 58000
$ws.Cells.Item(46, 10).Value = 58000
$ws.Cells.Item(46, 12).Value = 58000
$ws.Cells.Item(46, 14).Value = -58462
$ws.Cells.Item(54, 8).Value = 40025
$ws.Cells.Item(54, 9).Value = 29999
$ws.Cells.Item(54, 10).Value = 60077
$ws.Cells.Item(54, 11).Value = 29999
$ws.Cells.Item(54, 12).Value = 60077
$ws.Cells.Item(54, 13).Value = -29479
$ws.Cells.Item(54, 14).Value = -61117
$ws.Cells.Item(61, 8).Value = 3999
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 3999
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 3999
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(61, 14).Value = -4583
$ws.Cells.Item(75, 8).Value = 44995.5
$ws.Cells.Item(75, 9).Value = 39994
$ws.Cells.Item(75, 11).Value = 39994
$ws.Cells.Item(75, 13).Value = -39058
$ws.Cells.Item(78, 8).Value = 44995.5
$ws.Cells.Item(78, 9).Value = 39994
$ws.Cells.Item(78, 11).Value = 119982
$ws.Cells.Item(78, 13).Value = -115302
$ws.Cells.Item(96, 8).Value = 1462.2858
$ws.Cells.Item(96, 9).Value = 1462.2858
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 1462.2858
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 13).Value = -89.28580000000011
$ws.Cells.Item(96, 14).ClearContents()
$ws.Cells.Item(100, 8).Value = 100001970
$ws.Cells.Item(100, 9).Value = 1345.1428
$ws.Cells.Item(100, 10).Value = 333336740
$ws.Cells.Item(100, 11).Value = 2690.2856
$ws.Cells.Item(100, 12).Value = 666673480
$ws.Cells.Item(100, 13).Value = -2149.2856
$ws.Cells.Item(100, 14).Value = -666674562
$ws.Cells.Item(110, 8).Value = 49000.332
$ws.Cells.Item(110, 10).Value = 49000.332
$ws.Cells.Item(110, 12).Value = 49000.332
$ws.Cells.Item(110, 14).Value = -57180.332
$ws.Cells.Item(113, 8).Value = 932.625
$ws.Cells.Item(113, 9).Value = 1163.9166
$ws.Cells.Item(113, 10).Value = 238.75
$ws.Cells.Item(113, 11).Value = 3491.7498
$ws.Cells.Item(113, 12).Value = 716.25
$ws.Cells.Item(113, 13).Value = -1321.7498
$ws.Cells.Item(113, 14).Value = -5056.25
$ws.Cells.Item(132, 8).Value = 2031.0741
$ws.Cells.Item(132, 9).Value = 1899.7755
$ws.Cells.Item(132, 11).Value = 5699.3265
$ws.Cells.Item(132, 13).Value = -3169.3265
$ws.Cells.Item(134, 8).Value = 58000
$ws.Cells.Item(134, 10).Value = 58000
$ws.Cells.Item(134, 12).Value = 174000
$ws.Cells.Item(134, 14).Value = -179070
